# Get export date time from filename
#
# Rows A2:C3 on the "SomeMO1" sheet previously held placeholder text ("A")
# pulled from the shared-string table. They now hold actual numeric sample
# data (1,2,3 / 4,5,6), formatted like the existing numeric rows (6-7) with
# a Text number format, and that sheet becomes the active/selected sheet
# and range.

$wb = $excel.ActiveWorkbook

$wsSomeMO1 = $wb.Worksheets.Item("SomeMO1")

# Replace the placeholder "A" values with real numeric data.
$wsSomeMO1.Range("A2").Value = 1
$wsSomeMO1.Range("B2").Value = 2
$wsSomeMO1.Range("C2").Value = 3
$wsSomeMO1.Range("A3").Value = 4
$wsSomeMO1.Range("B3").Value = 5
$wsSomeMO1.Range("C3").Value = 6

# Match the number format already used by the numeric rows below (6-7).
$wsSomeMO1.Range("A2:C3").NumberFormat = "@"

# Make "SomeMO1" the active sheet/tab, with A2:C3 selected.
$wsSomeMO1.Activate()
$wsSomeMO1.Range("A2:C3").Select()
